$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

$ws.Range("D2").Value = "57.694.86"
$ws.Range("E2").Value = "  -1.21%  "

$ws.Range("D3").Value = "2.443.81"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").Value = "  -0.18%  "

Set-CellText "D5" "522.61"
$ws.Range("E5").Value = "  -0.01%  "

Set-CellText "D6" "129.75"
$ws.Range("E6").Value = "  -2.39%  "

$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("E8").Value = "  +0.39%  "

$ws.Range("B9").Value = "LidoStakedEther"
$ws.Range("C9").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D9").Value = "2.442.57"
$ws.Range("E9").Value = "  -3.21%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-CellText "D10" "0.0973"
$ws.Range("E10").Value = "  -0.33%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-CellText "D11" "0.151"
$ws.Range("E11").Value = "  -3.45%  "

$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText "D12" "4.89"
$ws.Range("E12").Value = "  -5.34%  "

$ws.Range("B13").Value = "Cardano"
$ws.Range("C13").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-CellText "D13" "0.321"
$ws.Range("E13").Value = "  -3.58%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.879.40"
$ws.Range("E14").Value = "  -2.03%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "57.640.68"
$ws.Range("E15").Value = "  -1.61%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-CellText "D16" "21.67"
$ws.Range("E16").Value = "  -2.30%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-CellText "D17" "0.0000132"
$ws.Range("E17").Value = "  -2.14%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.442.98"
$ws.Range("E18").Value = "  -2.88%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-CellText "D19" "10.32"
$ws.Range("E19").Value = "  -3.33%  "

$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText "D20" "4.12"
$ws.Range("E20").Value = "  -1.03%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-CellText "D21" "313.99"
$ws.Range("E21").Value = "  -2.64%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-CellText "D22" "6.08"
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-CellText "D23" "1.00"
$ws.Range("E23").Value = "  +0.40%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-CellText "D24" "64.90"
$ws.Range("E24").Value = "  +0.43%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-CellText "D25" "0.411"
$ws.Range("E25").Value = "  +1.83%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-CellText "D26" "1.00"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-CellText "D27" "0.156"
$ws.Range("E27").Value = "  -3.06%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-CellText "D28" "7.19"
$ws.Range("E28").Value = "  -3.00%  "

$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-CellText "D29" "172.27"
$ws.Range("E29").Value = "  +2.68%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0732"
$ws.Range("E30").Value = "  -3.13%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-CellText "D31" "1.69"
$ws.Range("E31").Value = "  -1.55%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-CellText "D32" "6.08"
$ws.Range("E32").Value = "  -2.98%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-CellText "D33" "1.14"
$ws.Range("E33").Value = "  -5.56%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-CellText "D34" "0.999"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-CellText "D35" "0.999"
$ws.Range("E35").Value = "  +0.61%  "

$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-CellText "D36" "17.78"
$ws.Range("E36").Value = "  -1.89%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-CellText "D37" "1.17"
$ws.Range("E37").Value = "  -6.72%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-CellText "D38" "3.78"
$ws.Range("E38").Value = "  -4.82%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-CellText "D39" "36.29"
$ws.Range("E39").Value = "  +0.43%  "

$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-CellText "D40" "1.45"
$ws.Range("E40").Value = "  -1.57%  "

$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-CellText "D41" "0.786"
$ws.Range("E41").Value = "  +0.78%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellText "D42" "3.39"
$ws.Range("E42").Value = "  -3.17%  "

$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-CellText "D43" "263.26"
$ws.Range("E43").Value = "  -5.76%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-CellText "D44" "0.581"
$ws.Range("E44").Value = "  -2.75%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText "D45" "4.78"
$ws.Range("E45").Value = "  -5.40%  "

$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-CellText "D47" "122.68"
$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-CellText "D48" "0.0491"
$ws.Range("E48").Value = "  -2.19%  "

$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText "D49" "0.0210"
$ws.Range("E49").Value = "  -1.77%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText "D50" "17.00"
$ws.Range("E50").Value = "  -4.52%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-CellText "D51" "16.24"
$ws.Range("E51").Value = "  -4.12%  "
